$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update B8/B9: these rows previously represented extr1/extr2 but now
# represent the newly introduced line7/line8 entries (two new shared
# strings are inserted into the table right after line6, shifting the
# extr* entries down). We also re-assert B10-B15 so their text stays
# correct (extr1..extr6) even though the underlying shared-string table
# is being rewritten.
$ws.Range("B8").Value2 = "line7"
$ws.Range("B9").Value2 = "line8"
$ws.Range("B10").Value2 = "extr1"
$ws.Range("B11").Value2 = "extr2"
$ws.Range("B12").Value2 = "extr3"
$ws.Range("B13").Value2 = "extr4"
$ws.Range("B14").Value2 = "extr5"
$ws.Range("B15").Value2 = "extr6"

# --- Update C/D/E values for rows 8-15 (existing rows) ---
$ws.Range("C8").Value2 = 14
$ws.Range("D8").Value2 = 11
$ws.Range("E8").Value2 = $true

$ws.Range("C9").Value2 = 16
$ws.Range("D9").Value2 = 9
$ws.Range("E9").Value2 = $true

$ws.Range("C10").Value2 = 5
$ws.Range("D10").Value2 = 12
$ws.Range("E10").Value2 = $true

$ws.Range("C11").Value2 = 5
$ws.Range("D11").Value2 = 9
$ws.Range("E11").Value2 = $true

$ws.Range("C12").Value2 = 10
$ws.Range("D12").Value2 = 11
$ws.Range("E12").Value2 = $false

$ws.Range("C13").Value2 = 7
$ws.Range("D13").Value2 = 8
$ws.Range("E13").Value2 = $false

$ws.Range("C14").Value2 = 9
$ws.Range("D14").Value2 = 11
$ws.Range("E14").Value2 = $true

$ws.Range("C15").Value2 = 7
$ws.Range("D15").Value2 = 11
$ws.Range("E15").Value2 = $true

# --- Append new rows 16 and 17 (extr7, extr8) ---
$ws.Range("A16").Value2 = 14
$ws.Range("B16").Value2 = "extr7"
$ws.Range("C16").Value2 = 5
$ws.Range("D16").Value2 = 7
$ws.Range("E16").Value2 = $false

$ws.Range("A17").Value2 = 15
$ws.Range("B17").Value2 = "extr8"
$ws.Range("C17").Value2 = 8
$ws.Range("D17").Value2 = 5
$ws.Range("E17").Value2 = $true

# Apply the same cell formatting (bold, centered, thin border) that the
# other column-A "index" cells use, by copying formats from A15.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
